$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 12200
$ws.Range("J54").Value = 22400
$ws.Range("L54").Value = 22400
$ws.Range("N54").Value = -23372
$ws.Range("H70").Value = 3734.5334
$ws.Range("I70").Value = 3159.5
$ws.Range("J70").Value = 4884.6
$ws.Range("K70").Value = 9478.5
$ws.Range("L70").Value = 14653.8
$ws.Range("M70").Value = -9208.5
$ws.Range("N70").Value = -15193.8
$ws.Range("H73").Value = 3734.5334
$ws.Range("I73").Value = 3159.5
$ws.Range("J73").Value = 4884.6
$ws.Range("K73").Value = 9478.5
$ws.Range("L73").Value = 14653.8
$ws.Range("M73").Value = -8542.5
$ws.Range("N73").Value = -16525.8
$ws.Range("H107").Value = 355.125
$ws.Range("I107").Value = 355.125
$ws.Range("K107").Value = 355.125
$ws.Range("M107").Value = 1564.875
$ws.Range("H111").Value = 5420.857
$ws.Range("J111").Value = 1966.3334
$ws.Range("L111").Value = 5899.0002
$ws.Range("N111").Value = -12033.0002
$ws.Range("H137").Value = 864457.1
$ws.Range("I137").Value = 1001923.56
$ws.Range("K137").Value = 3005770.68
$ws.Range("M137").Value = -3003220.68
$ws.Range("H138").Value = 3057.9592
$ws.Range("I138").Value = 1219.3
$ws.Range("J138").Value = 3529.4102
$ws.Range("K138").Value = 3657.9
$ws.Range("L138").Value = 10588.2306
$ws.Range("M138").Value = 1482.1
$ws.Range("N138").Value = -20868.2306
$ws.Range("H141").Value = 2229.1177
$ws.Range("I141").Value = 2014.1666
$ws.Range("J141").Value = 2745
$ws.Range("K141").Value = 6042.4998
$ws.Range("L141").Value = 8235
$ws.Range("M141").Value = -862.4997999999996
$ws.Range("N141").Value = -18595

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1743.1818
$ws.Range("I2").Value = 1575.2222
$ws.Range("J2").Value = 2499
$ws.Range("K2").Value = 1575.2222
$ws.Range("L2").Value = 2499
$ws.Range("M2").Value = -1462.2222
$ws.Range("N2").Value = -2725
$ws.Range("H32").Value = 2426.6099
$ws.Range("I32").Value = 2515.1282
$ws.Range("J32").Value = 700.5
$ws.Range("K32").Value = 2515.1282
$ws.Range("L32").Value = 700.5
$ws.Range("M32").Value = -2228.1282
$ws.Range("N32").Value = -1274.5
$ws.Range("H102").Value = 2412.476
$ws.Range("I102").Value = 2192.611
$ws.Range("K102").Value = 2192.611
$ws.Range("M102").Value = -570.6109999999999
$ws.Range("H116").Value = 1743.1818
$ws.Range("I116").Value = 1575.2222
$ws.Range("J116").Value = 2499
$ws.Range("K116").Value = 1575.2222
$ws.Range("L116").Value = 2499
$ws.Range("M116").Value = 718.7778000000001
$ws.Range("N116").Value = -7087
$ws.Range("H132").Value = 1797.6
$ws.Range("I132").Value = 1489.7142
$ws.Range("K132").Value = 4469.142599999999
$ws.Range("M132").Value = -1939.142599999999
$ws.Range("H139").Value = 76602.62
$ws.Range("I139").Value = 69998
$ws.Range("J139").Value = 77153
$ws.Range("K139").Value = 69998
$ws.Range("L139").Value = 77153
$ws.Range("M139").Value = -64858
$ws.Range("N139").Value = -87433
$ws.Range("H141").Value = 62333.332
$ws.Range("J141").Value = 62333.332
$ws.Range("L141").Value = 62333.332
$ws.Range("N141").Value = -72693.33199999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1743.1818
$ws.Range("I3").Value = 1575.2222
$ws.Range("J3").Value = 2499
$ws.Range("K3").Value = 1575.2222
$ws.Range("L3").Value = 2499
$ws.Range("M3").Value = -1461.2222
$ws.Range("N3").Value = -2727
$ws.Range("H20").Value = 25005458
$ws.Range("I20").Value = 45462016
$ws.Range("J20").Value = 2999.889
$ws.Range("K20").Value = 45462016
$ws.Range("L20").Value = 2999.889
$ws.Range("M20").Value = -45461769
$ws.Range("N20").Value = -3493.889
$ws.Range("H134").Value = 3082.6453
$ws.Range("I134").Value = 2873.96
$ws.Range("K134").Value = 8621.880000000001
$ws.Range("M134").Value = -6086.880000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 773.3333
$ws.Range("J4").Value = 1060
$ws.Range("L4").Value = 1060
$ws.Range("N4").Value = -1284
$ws.Range("H10").Value = 782.5
$ws.Range("I10").Value = 782.5
$ws.Range("K10").Value = 782.5
$ws.Range("M10").Value = -643.5
$ws.Range("H16").Value = 2105
$ws.Range("J16").Value = 1868.2858
$ws.Range("L16").Value = 1868.2858
$ws.Range("N16").Value = -2442.2858
$ws.Range("H51").Value = 59099
$ws.Range("J51").Value = 59099
$ws.Range("L51").Value = 59099
$ws.Range("N51").Value = -60571
$ws.Range("H61").Value = 59099
$ws.Range("J61").Value = 59099
$ws.Range("L61").Value = 59099
$ws.Range("N61").Value = -59795
$ws.Range("H62").Value = 12503523
$ws.Range("I62").Value = 12503523
$ws.Range("K62").Value = 12503523
$ws.Range("M62").Value = -12502899
$ws.Range("H65").Value = 12503523
$ws.Range("I65").Value = 12503523
$ws.Range("K65").Value = 62517615
$ws.Range("M65").Value = -62514495
$ws.Range("H113").Value = 2105
$ws.Range("J113").Value = 1868.2858
$ws.Range("L113").Value = 1868.2858
$ws.Range("N113").Value = -6208.2858
$ws.Range("H132").Value = 11118380
$ws.Range("I132").Value = 6418.5835
$ws.Range("J132").Value = 55566224
$ws.Range("K132").Value = 19255.7505
$ws.Range("L132").Value = 166698672
$ws.Range("M132").Value = -16725.7505
$ws.Range("N132").Value = -166703732
$ws.Range("H134").Value = 2668.4915
$ws.Range("I134").Value = 2308.551
$ws.Range("K134").Value = 6925.653
$ws.Range("M134").Value = -4390.653

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 164.25
$ws.Range("J12").Value = 187.4
$ws.Range("L12").Value = 562.2
$ws.Range("N12").Value = -908.2
$ws.Range("H22").Value = 6927.7144
$ws.Range("I22").Value = 5749.5
$ws.Range("K22").Value = 17248.5
$ws.Range("M22").Value = -17079.5
$ws.Range("H27").Value = 6927.7144
$ws.Range("I27").Value = 5749.5
$ws.Range("K27").Value = 17248.5
$ws.Range("M27").Value = -17146.5
$ws.Range("H54").Value = 5594.8
$ws.Range("I54").Value = 2984
$ws.Range("J54").Value = 6247.5
$ws.Range("K54").Value = 8952
$ws.Range("L54").Value = 18742.5
$ws.Range("M54").Value = -8393
$ws.Range("N54").Value = -19860.5
$ws.Range("H60").Value = 1179326.4
$ws.Range("I60").Value = 2857837.2
$ws.Range("K60").Value = 8573511.600000001
$ws.Range("M60").Value = -8573260.600000001
$ws.Range("H140").Value = 6742.7144
$ws.Range("I140").Value = 3722.8333
$ws.Range("J140").Value = 9007.625
$ws.Range("K140").Value = 11168.4999
$ws.Range("L140").Value = 27022.875
$ws.Range("M140").Value = -5988.499899999999
$ws.Range("N140").Value = -37382.875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 45000
$ws.Range("J63").Value = 45000
$ws.Range("L63").Value = 45000
$ws.Range("N63").Value = -46372
$ws.Range("H66").Value = 45000
$ws.Range("J66").Value = 45000
$ws.Range("L66").Value = 135000
$ws.Range("N66").Value = -141864
$ws.Range("H107").Value = 4791.1763
$ws.Range("I107").Value = 534.2857
$ws.Range("J107").Value = 7771
$ws.Range("K107").Value = 534.2857
$ws.Range("L107").Value = 7771
$ws.Range("M107").Value = 1385.7143
$ws.Range("N107").Value = -11611
$ws.Range("H113").Value = 7820.375
$ws.Range("I113").Value = 3832.3333
$ws.Range("J113").Value = 19784.5
$ws.Range("K113").Value = 3832.3333
$ws.Range("L113").Value = 19784.5
$ws.Range("M113").Value = -1662.3333
$ws.Range("N113").Value = -24124.5
$ws.Range("H132").Value = 2488.3333
$ws.Range("I132").Value = 2488.3333
$ws.Range("K132").Value = 7464.999899999999
$ws.Range("M132").Value = -4934.999899999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 30000
$ws.Range("I23").Value = 30000
$ws.Range("K23").Value = 30000
$ws.Range("M23").Value = -29770
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H46").Value = 4664.6665
$ws.Range("I46").Value = 5998.25
$ws.Range("J46").Value = 1997.5
$ws.Range("K46").Value = 5998.25
$ws.Range("L46").Value = 1997.5
$ws.Range("M46").Value = -5810.25
$ws.Range("N46").Value = -2373.5
$ws.Range("H61").Value = 3832.875
$ws.Range("I61").Value = 3824
$ws.Range("J61").Value = 3895
$ws.Range("K61").Value = 3824
$ws.Range("L61").Value = 3895
$ws.Range("M61").Value = -3622
$ws.Range("N61").Value = -4299
$ws.Range("H113").Value = 3832.875
$ws.Range("I113").Value = 3824
$ws.Range("J113").Value = 3895
$ws.Range("K113").Value = 3824
$ws.Range("L113").Value = 3895
$ws.Range("M113").Value = -1654
$ws.Range("N113").Value = -8235
$ws.Range("H122").Value = 12991
$ws.Range("I122").Value = 6975
$ws.Range("J122").Value = 15999
$ws.Range("K122").Value = 20925
$ws.Range("L122").Value = 47997
$ws.Range("M122").Value = -18475
$ws.Range("N122").Value = -52897
$ws.Range("H136").Value = 6209.5625
$ws.Range("I136").Value = 4279.4165
$ws.Range("K136").Value = 12838.2495
$ws.Range("M136").Value = -10288.2495

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4243.625
$ws.Range("I96").Value = 2278.4285
$ws.Range("K96").Value = 2278.4285
$ws.Range("M96").Value = -905.4285
$ws.Range("H122").Value = 19232570
$ws.Range("J122").Value = 83334536
$ws.Range("L122").Value = 250003608
$ws.Range("N122").Value = -250008508
$ws.Range("H126").Value = 1120
$ws.Range("I126").Value = 1300
$ws.Range("J126").Value = 400
$ws.Range("K126").Value = 3900
$ws.Range("L126").Value = 1200
$ws.Range("M126").Value = -1430
$ws.Range("N126").Value = -6140
